# Updates cryptos list values (price + 1h volume change) per commit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.715.60"
$ws.Range("E2").Value = "  -2.75%  "

$ws.Range("D3").Value = "1.743.80"
$ws.Range("E3").Value = "  -4.95%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9998"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.19%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "238.92"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -8.11%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9999"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.18%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5051"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -5.16%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "41.79"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -6.95%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2658"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -11.55%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06130"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -10.60%  "

$ws.Range("D11").Value = "1.745.59"
$ws.Range("E11").Value = "  -4.95%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.06955"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -4.34%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "15.27"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -12.68%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.495"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -9.22%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5975"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -18.62%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "76.62"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -13.73%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9994"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.34%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9996"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.20%  "

$ws.Range("D19").Value = "25.710.32"
$ws.Range("E19").Value = "  -2.86%  "

$ws.Range("E20").Value = "  -16.16%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.000006771"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -13.87%  "

$ws.Range("D22").Value = "1.964.83"
$ws.Range("E22").Value = "  -5.52%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.042"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -11.60%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.162"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -11.47%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.130"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -13.87%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "137.82"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -3.38%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.519"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -9.61%  "

$ws.Range("E28").Value = "  -16.79%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "14.98"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -11.38%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "103.26"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -6.39%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.754"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -10.86%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08105"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -7.88%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.456"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -13.52%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04499"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -5.89%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9987"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.21%  "

$ws.Range("E36").Value = "  -9.40%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9833"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -12.72%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.6108"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -16.13%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.657"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -13.94%  "

$ws.Range("E40").Value = "  -8.88%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.915"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -16.27%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9999"
$ws.Range("D42").ClearFormats()

$ws.Range("E43").Value = "  -3.48%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3801"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -19.14%  "

$ws.Range("E45").Value = "  -12.93%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.7282"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -19.42%  "


$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1112"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -9.32%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "30.11"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -13.22%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.889"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -19.71%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "52.46"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -12.50%  "
